$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row labels (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize Spanish connector words (de/del/el/la/los/las/y) to capitalized
#     form (De/Del/El/La/Los/Las/Y) within municipality (and a couple of state)
#     names, mirroring a capitalization-style cleanup pass over the data. ---
$changes = @(
    @(7, "B", "Pabellón De Arteaga"),
    @(8, "B", "Rincón De Romos"),
    @(9, "B", "San Francisco De Los Romo"),
    @(10, "B", "San José De Gracia"),
    @(31, "B", "Amatenango De La Frontera"),
    @(34, "B", "Bejucal De Ocampo"),
    @(41, "B", "Chiapa De Corzo"),
    @(47, "B", "Comitán De Domínguez"),
    @(64, "B", "Mazapa De Madero"),
    @(67, "B", "Montecristo De Guerrero"),
    @(70, "B", "Ocozocoautla De Espinosa"),
    @(79, "B", "Salto De Agua"),
    @(80, "B", "San Cristóbal De Las Casas"),
    @(114, "B", "Hidalgo Del Parral"),
    @(122, "B", "San Francisco De Borja"),
    @(124, "B", "Valle De Zaragoza"),
    @(140, "B", "San Juan De Sabinas"),
    @(152, "A", "Ciudad De México"),
    @(169, "B", "Coneto De Comonfort"),
    @(181, "B", "Nombre De Dios"),
    @(185, "B", "Pánuco De Coronado"),
    @(190, "B", "San Juan Del Río"),
    @(195, "A", "Estado De México"),
    @(195, "B", "Acambay De Ruíz Castañeda"),
    @(198, "B", "Almoloya De Alquisiras"),
    @(199, "B", "Almoloya De Juárez"),
    @(204, "B", "Atizapán De Zaragoza"),
    @(211, "B", "Coacalco De Berriozábal"),
    @(216, "B", "Ecatepec De Morelos"),
    @(220, "B", "Ixtapan De La Sal"),
    @(230, "B", "Naucalpan De Juárez"),
    @(237, "B", "San Felipe Del Progreso"),
    @(238, "B", "San Simón De Guerrero"),
    @(240, "B", "Soyaniquilpan De Juárez"),
    @(247, "B", "Tenango Del Valle"),
    @(257, "B", "Tlalnepantla De Baz"),
    @(262, "B", "Valle De Chalco Solidaridad"),
    @(263, "B", "Villa De Allende"),
    @(264, "B", "Villa Del Carbón"),
    @(273, "B", "San Miguel De Allende"),
    @(274, "B", "Apaseo El Alto"),
    @(275, "B", "Apaseo El Grande"),
    @(283, "B", "Dolores Hidalgo Cuna De La Independencia Nacional"),
    @(293, "B", "Purísima Del Rincón"),
    @(297, "B", "San Diego De La Unión"),
    @(299, "B", "San Francisco Del Rincón"),
    @(301, "B", "San Luis De La Paz"),
    @(302, "B", "Santa Cruz De Juventino Rosas"),
    @(303, "B", "Silao De La Victoria"),
    @(308, "B", "Valle De Santiago"),
    @(314, "B", "Acapulco De Juárez"),
    @(316, "B", "Ajuchitlán Del Progreso"),
    @(317, "B", "Alcozauca De Guerrero"),
    @(320, "B", "Atenango Del Río"),
    @(322, "B", "Atoyac De Álvarez"),
    @(323, "B", "Ayutla De Los Libres"),
    @(326, "B", "Chilapa De Álvarez"),
    @(327, "B", "Chilpancingo De Los Bravo"),
    @(331, "B", "Coyuca De Benítez"),
    @(332, "B", "Coyuca De Catalán"),
    @(335, "B", "Cutzamala De Pinzón"),
    @(341, "B", "Huitzuco De Los Figueroa"),
    @(342, "B", "Iguala De La Independencia"),
    @(343, "B", "Ixcateopan De Cuauhtémoc"),
    @(344, "B", "Zihuatanejo De Azueta"),
    @(346, "B", "La Unión De Isidoro Montes De Oca"),
    @(349, "B", "Mártir De Cuilapan"),
    @(361, "B", "Taxco De Alarcón"),
    @(363, "B", "Técpan De Galeana"),
    @(365, "B", "Tepecoacuilco De Trujano"),
    @(367, "B", "Tixtla De Guerrero"),
    @(371, "B", "Tlapa De Comonfort"),
    @(385, "B", "Atotonilco De Tula"),
    @(386, "B", "Atotonilco El Grande"),
    @(391, "B", "Cuautepec De Hinojosa"),
    @(395, "B", "Huasca De Ocampo"),
    @(398, "B", "Huejutla De Reyes"),
    @(401, "B", "Jacala De Ledezma"),
    @(407, "B", "Mineral Del Chico"),
    @(408, "B", "Mineral Del Monte"),
    @(409, "B", "Mixquiahuala De Juárez"),
    @(410, "B", "Molango De Escamilla"),
    @(412, "B", "Nopala De Villagrán"),
    @(413, "B", "Omitlán De Juárez"),
    @(414, "B", "Pachuca De Soto"),
    @(417, "B", "Progreso De Obregón"),
    @(421, "B", "Santiago De Anaya"),
    @(425, "B", "Tenango De Doria"),
    @(427, "B", "Tepehuacán De Guerrero"),
    @(428, "B", "Tepeji Del Río De Ocampo"),
    @(429, "B", "Tezontepec De Aldama"),
    @(434, "B", "Tula De Allende"),
    @(435, "B", "Tulancingo De Bravo"),
    @(437, "B", "Zacualtipán De Ángeles"),
    @(442, "B", "Ahualulco De Mercado"),
    @(446, "B", "Atotonilco El Alto"),
    @(448, "B", "Autlán De Navarro"),
    @(451, "B", "Cañadas De Obregón"),
    @(456, "B", "Concepción De Buenos Aires"),
    @(461, "B", "Encarnación De Díaz"),
    @(464, "B", "Huejuquilla El Alto"),
    @(471, "B", "Lagos De Moreno"),
    @(475, "B", "Ojuelos De Jalisco"),
    @(480, "B", "San Juan De Los Lagos"),
    @(483, "B", "San Miguel El Alto"),
    @(484, "B", "Santa María De Los Ángeles"),
    @(489, "B", "Tepatitlán De Morelos"),
    @(491, "B", "Tizapán El Alto"),
    @(492, "B", "Tlajomulco De Zúñiga"),
    @(501, "B", "Unión De San Antonio"),
    @(502, "B", "Valle De Guadalupe"),
    @(505, "B", "Yahualica De González Gallo"),
    @(508, "B", "Zapotitlán De Vadillo"),
    @(509, "B", "Zapotlán Del Rey"),
    @(510, "B", "Zapotlán El Grande"),
    @(583, "B", "Tiquicheo De Nicolás Romero"),
    @(604, "B", "Coatlán Del Río"),
    @(613, "B", "Puente De Ixtla"),
    @(617, "B", "Tetela Del Volcán"),
    @(618, "B", "Tlaltizapán De Zapata"),
    @(628, "B", "Amatlán De Cañas"),
    @(630, "B", "Ixtlán Del Río"),
    @(636, "B", "Santa María Del Oro"),
    @(644, "B", "Ciénega De Flores"),
    @(657, "B", "San Nicolás De Los Garza"),
    @(660, "B", "Acatlán De Pérez Figueroa"),
    @(665, "B", "Chalcatongo De Hidalgo"),
    @(668, "B", "Coicoyán De Las Flores"),
    @(671, "B", "Guadalupe De Ramírez"),
    @(672, "B", "Guevea De Humboldt"),
    @(673, "B", "Heroica Ciudad De Ejutla De Crespo"),
    @(674, "B", "Heroica Ciudad De Huajuapan De León"),
    @(675, "B", "Heroica Ciudad De Tlaxiaco"),
    @(677, "B", "Ixtlán De Juárez"),
    @(678, "B", "Heroica Ciudad De Juchitán De Zaragoza"),
    @(682, "B", "Mártires De Tacubaya"),
    @(684, "B", "Mazatlán Villa De Flores"),
    @(685, "B", "Miahuatlán De Porfirio Díaz"),
    @(687, "B", "Nejapa De Madero"),
    @(688, "B", "Oaxaca De Juárez"),
    @(689, "B", "Ocotlán De Morelos"),
    @(690, "B", "Pinotepa De Don Luis"),
    @(692, "B", "Putla Villa De Guerrero"),
    @(693, "B", "Reforma De Pineda"),
    @(699, "B", "San Antonino El Alto"),
    @(700, "B", "San Antonio De La Cal"),
    @(703, "B", "San Baltazar Yatzachi El Bajo"),
    @(710, "B", "San Felipe Jalapa De Díaz"),
    @(722, "B", "San Juan Bautista Lo De Soto"),
    @(749, "B", "San Miguel Del Puerto"),
    @(750, "B", "San Miguel El Grande"),
    @(758, "B", "San Pablo Villa De Mitla"),
    @(760, "B", "San Pedro El Alto"),
    @(770, "B", "San Pedro Y San Pablo Teposcolula"),
    @(790, "B", "Santa María Del Tule"),
    @(796, "B", "Santa María Jalapa Del Marqués"),
    @(821, "B", "Santo Domingo De Morelos"),
    @(833, "B", "Tamazulápam Del Espíritu Santo"),
    @(834, "B", "Tataltepec De Valdés"),
    @(835, "B", "Tepelmeme Villa De Morelos"),
    @(836, "B", "Tlacolula De Matamoros"),
    @(838, "B", "Villa De Etla"),
    @(839, "B", "Villa De Tututepec De Melchor Ocampo"),
    @(840, "B", "Villa Sola De Vega"),
    @(842, "B", "Zapotitlán Del Río"),
    @(844, "B", "Zimatlán De Álvarez"),
    @(854, "B", "Chalchicomula De Sesma"),
    @(869, "B", "Huehuetlán El Chico"),
    @(873, "B", "Izúcar De Matamoros"),
    @(878, "B", "Los Reyes De Juárez"),
    @(883, "B", "Palmar De Bravo"),
    @(894, "B", "San Salvador El Seco"),
    @(895, "B", "San Salvador El Verde"),
    @(899, "B", "Tecali De Herrera"),
    @(904, "B", "Tepanco De López"),
    @(907, "B", "Teteles De Avila Castillo"),
    @(910, "B", "Tlacotepec De Benito Juárez"),
    @(916, "B", "Xayacatlán De Bravo"),
    @(927, "B", "Amealco De Bonfil"),
    @(929, "B", "Cadereyta De Montes"),
    @(934, "B", "Jalpan De Serra"),
    @(935, "B", "Landa De Matamoros"),
    @(937, "B", "Pinal De Amoles"),
    @(940, "B", "San Juan Del Río"),
    @(951, "B", "Axtla De Terrazas"),
    @(957, "B", "Ciudad Del Maíz"),
    @(966, "B", "Mexquitic De Carmona"),
    @(972, "B", "San Ciro De Acosta"),
    @(978, "B", "Santa María Del Río"),
    @(980, "B", "Soledad De Graciano Sánchez"),
    @(989, "B", "Villa De Arista"),
    @(990, "B", "Villa De Arriaga"),
    @(991, "B", "Villa De Guadalupe"),
    @(992, "B", "Villa De La Paz"),
    @(993, "B", "Villa De Ramos"),
    @(994, "B", "Villa De Reyes"),
    @(1017, "B", "Nacozari De García"),
    @(1029, "B", "Jalpa De Méndez"),
    @(1059, "B", "Soto La Marina"),
    @(1073, "B", "Ixtacuixtla De Mariano Matamoros"),
    @(1074, "B", "Muñoz De Domingo Arenas"),
    @(1089, "B", "Alto Lucero De Gutiérrez Barrios"),
    @(1092, "B", "Amatlán De Los Reyes"),
    @(1100, "B", "Boca Del Río"),
    @(1102, "B", "Camarón De Tejeda"),
    @(1107, "B", "Cazones De Herrera"),
    @(1119, "B", "Cosamaloapan De Carpio"),
    @(1135, "B", "Hueyapan De Ocampo"),
    @(1136, "B", "Ignacio De La Llave"),
    @(1139, "B", "Ixhuatlán Del Café"),
    @(1146, "B", "Juchique De Ferrer"),
    @(1150, "B", "Las Vigas De Ramírez"),
    @(1153, "B", "Martínez De La Torre"),
    @(1155, "B", "Medellín De Bravo"),
    @(1159, "B", "Nanchital De Lázaro Cárdenas Del Río"),
    @(1164, "B", "Ozuluama De Mascareñas"),
    @(1168, "B", "Paso Del Macho"),
    @(1172, "B", "Poza Rica De Hidalgo"),
    @(1178, "B", "Sayula De Alemán"),
    @(1181, "B", "Soledad De Doblado"),
    @(1208, "B", "Vega De Alatorre"),
    @(1229, "B", "El Plateado De Joaquín Amaro"),
    @(1246, "B", "Moyahua De Estrada"),
    @(1247, "B", "Nochistlán De Mejía"),
    @(1248, "B", "Noria De Ángeles"),
    @(1256, "B", "Teúl De González Ortega"),
    @(1257, "B", "Tlaltenango De Sánchez Román"),
    @(1259, "B", "Villa De Cos"),
)

foreach ($change in $changes) {
    $row = $change[0]
    $col = $change[1]
    $val = $change[2]
    $ws.Range("$col$row").Value = $val
}

# --- Tiny floating point precision corrections on the two "Total" percentage
#     rows (last-bit recomputation artifacts) ---
$ws.Range("D313").Value = 0.09633911368015416
$ws.Range("D1218").Value = 0.09479768786127168

# --- Remove trailing footnote/metadata rows 1268:1272, shrinking the used
#     range down to A1:D1266 ---
$ws.Rows.Item(1268).Resize(5).Delete() | Out-Null
